$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.990.28"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").Value = "3.520.14"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.10"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.98"
$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("D7").Value = "3.517.28"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.84"
$ws.Range("E11").Value = "  +4.30%  "

$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").Value = "4.116.46"
$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.63"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").Value = "3.517.92"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "66.989.92"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  +9.21%  "

$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "437.16"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  -2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.71"
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("D25").Value = "3.660.17"
$ws.Range("E25").Value = "  +0.83%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.37"
$ws.Range("E29").Value = "  -3.15%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.39"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").Value = "3.513.92"
$ws.Range("E35").Value = "  +0.97%  "

$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.91"
$ws.Range("E37").Value = "  -3.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.04"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0891"
$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.91"
$ws.Range("E42").Value = "  -2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("E44").Value = "  -9.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.896"
$ws.Range("E45").Value = "  +1.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.05"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.32"
$ws.Range("E47").Value = "  +2.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.24"
$ws.Range("E48").Value = "  -5.32%  "

$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("E50").Value = "  -2.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.991"
$ws.Range("E51").Value = "  +0.39%  "

